$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Hjemme passive updated meanEMG legmaxROM": refresh the first four subject
# columns (B:E) on all three rows (header ids + the two condition rows,
# CON/STR) with the re-exported values for subjects 15 and 16.
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = -6.38301442675305
$ws.Range("C2").Value = 8.4136337932082483
$ws.Range("D2").Value = 14.454946235178909
$ws.Range("E2").Value = 23.730625264933053

$ws.Range("B3").Value = -23.802997961765413
$ws.Range("C3").Value = 11.526169832176834
$ws.Range("D3").Value = 34.995288937225325
$ws.Range("E3").Value = 12.154062769443755

$ws.Range("B1:E3").Select()
